# Eksperimenti za casovno zahtevnost
# Applies the changes described by the commit: adds a second experiment
# table (kmax-based) below the existing one, a matching "Grafikon 2"
# scatter chart with a linear trendline, a numeric display format for
# a handful of the new time values, and small title/label tweaks to the
# the existing chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the stray trailing "3000" row (B27) left over from the old
#    table - the row disappears entirely once its only cell is cleared.
# ---------------------------------------------------------------------
$ws.Range("B27").ClearContents()

# ---------------------------------------------------------------------
# 2. New kmax-vs-time experiment table, rows 32-58.
# ---------------------------------------------------------------------

# Header row + the little kmax/n/odstranjenih/izbranih summary block
$ws.Range("B32").Value = "kmax"
$ws.Range("C32").Value = "Čas (s)"
$ws.Range("F32").Value = "grafov "
$ws.Range("G32").Value = 1
$ws.Range("F33").Value = "n"
$ws.Range("G33").Value = 100
$ws.Range("F34").Value = "št. dreves za odstranitev"
$ws.Range("G34").Value = 48
$ws.Range("F35").Value = "št. dreves za izbiro"
$ws.Range("G35").Value = 50

# kmax (col B) / time-in-seconds (col C) pairs
$kmax = @(1,2,3,4,5,6,7,8,9,10,15,20,25,30,35,40,50,60,75,100,150,200,300,400,500,750)
$times = @(1.845,1.8580000000000001,1.91,2.0249999999999999,2.1110000000000002,2.1520000000000001,2.2240000000000002,2.4089999999999998,2.41,2.524,2.8079999999999998,3.11,3.4140000000000001,3.7120000000000002,4.3730000000000002,4.4349999999999996,5.032,5.5279999999999996,6.7069999999999999,8.2230000000000008,11.409000000000001,15.042,22.312999999999999,30.414000000000001,37.070999999999998,52.57)

# Rows that carry the explicit "0.000" display format in the source file
$formattedRows = @(33,35,39,41,43,45,50)

for ($i = 0; $i -lt $kmax.Length; $i++) {
    $row = 33 + $i
    $ws.Cells.Item($row, 2).Value = $kmax[$i]
    $ws.Cells.Item($row, 3).Value = $times[$i]
    if ($formattedRows -contains $row) {
        $ws.Cells.Item($row, 3).NumberFormat = "0.000"
    }
}

# ---------------------------------------------------------------------
# 3. Tidy up the selection / scroll position left by the editing author.
# ---------------------------------------------------------------------
$ws.Range("G4").Select()

# ---------------------------------------------------------------------
# 4. Existing chart ("Chart 1"): give it an explicit title.
# ---------------------------------------------------------------------
$chart1 = $ws.ChartObjects(1).Chart
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Časovna zahtevnost P1 v odvisnosti od n"

# ---------------------------------------------------------------------
# 5. New chart ("Grafikon 2"): kmax (x) vs. time (y) scatter plot with a
#    linear trendline, built from the table just written above.
# ---------------------------------------------------------------------
$co2 = $ws.ChartObjects().Add(400, 20, 300, 250)
$co2.Name = "Grafikon 2"

$chart2 = $co2.Chart
$chart2.ChartType = 74  # xlXYScatterLines

$series2 = $chart2.SeriesCollection().NewSeries()
$series2.Name = "=Sheet1!`$C`$32"
$series2.XValues = $ws.Range("B33:B58")
$series2.Values = $ws.Range("C33:C58")

$series2.Trendlines().Add(1) | Out-Null

$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Časovna zahtevnost v P1 v odvisnosti od kmax"
